function Set-RowValues($ws, $r, $a, $b, $c, $d, $e) {
    $ws.Cells.Item($r,1).Value = $a
    $ws.Cells.Item($r,2).Value = $b
    $ws.Cells.Item($r,3).Value = $c
    $ws.Cells.Item($r,4).Value = $d
    $ws.Cells.Item($r,5).Value = $e
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "LP1912" (sheet 1) - main schedule table
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

# Header updates
$ws1.Cells.Item(2,1).Value = "Última actualización: 07:24:45"
$ws1.Cells.Item(3,1).Value = "Total filas: 70"

# Insert the 5 new single-row entries scattered through the table, from the
# top down so every later insertion point is expressed in "already shifted"
# row numbers.

# 1) New row at row 48
$ws1.Rows.Item(48).Insert()
Set-RowValues $ws1 48 "07:24:45" "07:31" "16_SANTA ANA" 7 "LP1912"

# 2) New row at row 51 (was row 50 before the previous insert)
$ws1.Rows.Item(51).Insert()
Set-RowValues $ws1 51 "07:24:45" "07:34" "23_HERNANDEZ" 10 "LP1912"

# 3) New row at row 54 (was row 52 before the previous two inserts)
$ws1.Rows.Item(54).Insert()
Set-RowValues $ws1 54 "07:24:45" "07:46" "14_ABASTO" 22 "LP1912"

# 4) New row at row 58 (was row 55 before the previous three inserts)
$ws1.Rows.Item(58).Insert()
Set-RowValues $ws1 58 "07:24:45" "08:03" "11_ETCHEVERRY" 39 "LP1912"

# 5) New row at row 69 (was row 65 before the previous four inserts)
$ws1.Rows.Item(69).Insert()
Set-RowValues $ws1 69 "07:24:45" "08:43" "14_ABASTO" 79 "LP1912"

# 6) Five brand-new rows appended after the current last row (71-75)
Set-RowValues $ws1 71 "07:24:45" "09:01" "215A_EL PATO" 97 "LP1912"
Set-RowValues $ws1 72 "07:24:45" "09:10" "16_P MOR-SANTA ANA" 106 "LP1912"
Set-RowValues $ws1 73 "07:24:45" "09:16" "27_EL RETIRO" 112 "LP1912"
Set-RowValues $ws1 74 "07:24:45" "09:21" "26_HERNANDEZ" 117 "LP1912"
Set-RowValues $ws1 75 "07:24:45" "09:22" "17_ROMERO" 118 "LP1912"

# ---------------------------------------------------------------------------
# Sheet "LP1912-215" (sheet 2) - filtered 215-line subset
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Cells.Item(2,1).Value = "Última actualización: 07:24:45"
$ws2.Cells.Item(3,1).Value = "Total filas: 14"

# One brand-new row appended at the end (row 19)
Set-RowValues $ws2 19 "07:24:45" "09:01" "215A_EL PATO" 97 "LP1912"

# ---------------------------------------------------------------------------
# Sheet "6203-6173" (sheet 3) - L6203 / L6173 lines
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Cells.Item(2,1).Value = "Última actualización: 07:24:45"
$ws3.Cells.Item(3,1).Value = "Total filas: 17"

# 1) New row at row 16
$ws3.Rows.Item(16).Insert()
Set-RowValues $ws3 16 "07:24:45" "07:38" "215A_LA PLATA" 14 "L6173"

# 2) New row at row 19 (was row 18 before the previous insert)
$ws3.Rows.Item(19).Insert()
Set-RowValues $ws3 19 "07:24:45" "08:35" "215A_LA PLATA" 71 "L6173"

# 3) Brand-new row appended at the end (row 22)
Set-RowValues $ws3 22 "07:24:45" "09:08" "215D_LA PLATA" 104 "L6203"
